# edit.ps1 — applies the authored changes to the presentation:
#   1. Refresh the cached "today" date fields (3/10/2025 -> 3/16/2025)
#      on every slide layout, the slide master, the handout master and
#      the notes master.
#   2. Split the run "Use bubble short algorithm." on slide 4 into
#      "Use bubble " / "sort " / "algorithm." (fixing short -> sort).
#   3. Remove every picture shape that was pasted onto slides 2-9.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Refresh cached date fields wherever a Date placeholder (type 16,
#    ppPlaceholderDate) is found, on the slide master, every custom
#    layout under it, the handout master and the notes master.
# ---------------------------------------------------------------------
$NEW_DATE = "3/16/2025"

function Update-DatePlaceholder($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -match "^\s*3/10/2025\s*$") {
                $tr.Text = $NEW_DATE
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li)
}

Update-DatePlaceholder $p.HandoutMaster
Update-DatePlaceholder $p.NotesMaster

# ---------------------------------------------------------------------
# 2. Slide 4: "Use bubble short algorithm." -> "Use bubble sort algorithm."
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
for ($i = 1; $i -le $s4.Shapes.Count; $i++) {
    $sh = $s4.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        for ($pi = 1; $pi -le $tr.Paragraphs().Count; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            $paraText = $para.Text.TrimEnd("`r", "`n")
            if ($paraText -eq "Use bubble short algorithm.") {
                $sub = $para.Characters(12, 6)
                $sub.Text = "sort "
            }
        }
    }
}

# ---------------------------------------------------------------------
# 3. Remove picture shapes added to slides 2-9 (msoPicture = 13).
# ---------------------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Type -eq 13) {
            $sh.Delete()
        }
    }
}
